$d = $word.ActiveDocument

# Locate the "Role 2:" paragraph (the last paragraph in the body before sectPr).
$role2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Role 2:") {
        $role2 = $p
    }
}

$lines = @(
    @{ Text = "TAGS:"; Bold = $true },
    @{ Text = "1 Taangaq -- Alcohol"; Bold = $false },
    @{ Text = "1 Alerquutet, Ayuqucirtuutet -- Instructions"; Bold = $false },
    @{ Text = "1 Qessaicaraq -- Not Being Lazy"; Bold = $false },
    @{ Text = "1 Tan'gaurluut Nasaurluut-llu Allakarluteng -- Boys and Girls Separate"; Bold = $false },
    @{ Text = "1 Piicak -- Prayer"; Bold = $false },
    @{ Text = "1 Ukverput, Agayuliyaraq -- Spirituality"; Bold = $false },
    @{ Text = "1 Qanruyutet, Qaneryarat -- Traditional Wisdom, Wise Words"; Bold = $false },
    @{ Text = "1 Anglicarillerkaq, Tukercaryaraq -- Child Rearing"; Bold = $false }
)

$cur = $role2
foreach ($line in $lines) {
    $cur.Range.InsertParagraphAfter()
    $cur = $d.Paragraphs.Last
    $nr = $cur.Range
    if ($line.Bold) {
        $nr.Font.Bold = 1
    }
    $nr.InsertBefore($line.Text)
}
